$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.775.14'
$ws.Range("E2").Value = '  +4.41%  '

$ws.Range("D3").Value = '3.337.46'
$ws.Range("E3").Value = '  +4.25%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.50'
$ws.Range("E5").Value = '  +3.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.85'
$ws.Range("E6").Value = '  +4.44%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +2.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.50'
$ws.Range("E9").Value = '  +2.35%  '

$ws.Range("E10").Value = '  +3.91%  '

$ws.Range("E11").Value = '  +1.53%  '

$ws.Range("D12").Value = '3.907.94'
$ws.Range("E12").Value = '  +4.21%  '

$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("E14").Value = '  +3.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.77'
$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '62.812.10'
$ws.Range("E16").Value = '  +4.45%  '

$ws.Range("D17").Value = '3.320.12'
$ws.Range("E17").Value = '  +3.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.48'
$ws.Range("E18").Value = '  +4.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.73'
$ws.Range("E19").Value = '  +4.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.46'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.89'
$ws.Range("E21").Value = '  +1.82%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.538'
$ws.Range("E23").Value = '  +1.49%  '

$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.181'
$ws.Range("E25").Value = '  +4.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.81'
$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("D27").Value = '0.0₃0959'
$ws.Range("E27").Value = '  +6.02%  '

$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.47'
$ws.Range("E29").Value = '  +4.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.97'
$ws.Range("E31").Value = '  +2.15%  '

$ws.Range("E32").Value = '  +2.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("E33").Value = '  +5.61%  '

$ws.Range("E34").Value = '  +2.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '160.57'
$ws.Range("E35").Value = '  +2.43%  '

$ws.Range("E36").Value = '  +9.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  +10.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.34'
$ws.Range("E38").Value = '  +6.37%  '

$ws.Range("D39").Value = '2.846.41'
$ws.Range("E39").Value = '  +2.56%  '

$ws.Range("E40").Value = '  +3.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0313'
$ws.Range("E41").Value = '  +8.55%  '

$ws.Range("E42").Value = '  +0.81%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.76'
$ws.Range("E43").Value = '  +2.59%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.749'
$ws.Range("E44").Value = '  +2.45%  '

$ws.Range("E45").Value = '  +2.55%  '

$ws.Range("D46").Value = '3.378.30'
$ws.Range("E46").Value = '  +4.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.96'
$ws.Range("E47").Value = '  +6.31%  '

$ws.Range("E48").Value = '  +3.32%  '

$ws.Range("E49").Value = '  +1.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.801'
$ws.Range("E50").Value = '  +0.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '282.76'
$ws.Range("E51").Value = '  +7.82%  '
